$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 1.78
    "H2"  = 3.5
    "I2"  = 4.3
    "J2"  = 2.32
    "L2"  = 4.65
    "M2"  = 1.07
    "N2"  = 6.8
    "O2"  = 1.35
    "P2"  = 2.95
    "Q2"  = 2.02
    "R2"  = 1.72
    "S2"  = 1.4
    "T2"  = 2.7
    "U2"  = 1.93
    "V2"  = 1.78
    "W2"  = 6.3
    "X2"  = 7.8
    "Z2"  = 14
    "AA2" = 15
    "AB2" = 30
    "AC2" = 6.8
    "AD2" = 6.8
    "AF2" = 90
    "AG2" = 800
    "AH2" = 10.75
    "AI2" = 23
    "AJ2" = 14.5
    "AK2" = 70
    "AL2" = 45
    "AN2" = 3.55
    "AO2" = 8.75
    "AP2" = 18.5
    "AQ2" = 30
    "AS2" = 250
    "AT2" = 2.7
    "AU2" = 7.7
    "AV2" = 75
    "AW2" = 6
    "AX2" = 25
    "AY2" = 32
    "AZ2" = 150
    "BB2" = 450
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
